# Update examples and code revision
# - cogen_model.xlsx: Processes sheet columns were re-ordered (description/type
#   columns swapped), a couple of column widths were tweaked, the
#   "cgam_processes" named range was narrowed from E6 to D6 and the active
#   sheet/tab moved from "Format" to "Processes".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Processes worksheet: swap the "description" (B) and "type" (E)
#    columns - both the header row and the six data rows.
# ---------------------------------------------------------------------
$processes = $wb.Worksheets.Item("Processes")

$processes.Range("B1").Value = "type"
$processes.Range("E1").Value = "description"

$processes.Range("B2").Value = "PRODUCTIVE"
$processes.Range("E2").Value = "Boiler"

$processes.Range("B3").Value = "PRODUCTIVE"
$processes.Range("E3").Value = "Turbine"

$processes.Range("B4").Value = "PRODUCTIVE"
$processes.Range("E4").Value = "Heat Exchanger"

$processes.Range("B5").Value = "PRODUCTIVE"
$processes.Range("E5").Value = "Pump"

$processes.Range("B6").Value = "PRODUCTIVE"
$processes.Range("E6").Value = "Alternator"

# Column width tweaks (B, E and G got resized).
$processes.Columns.Item(2).ColumnWidth = 13.333333333333334
$processes.Columns.Item(5).ColumnWidth = 14.666666666666666
$processes.Columns.Item(7).ColumnWidth = 13.0

# ---------------------------------------------------------------------
# 2. Named range "cgam_processes" shrinks from $A$1:$E$6 to $A$1:$D$6
#    now that the table is only 4 data columns wide (description/type
#    columns were merged/re-ordered).
# ---------------------------------------------------------------------
$names = $wb.Names
$cnt = $names.Count()
for ($i = 1; $i -le $cnt; $i++) {
  $n = $names.Item($i)
  $nm = $n.Name()
  if ($nm -eq "Processes!cgam_processes") {
    $n.RefersTo = "=Processes!`$A`$1:`$D`$6"
  }
}

# ---------------------------------------------------------------------
# 3. Active sheet / selection moved from "Format" to "Processes".
# ---------------------------------------------------------------------
$processes.Activate()
$processes.Range("D4").Select()

$format = $wb.Worksheets.Item("Format")
$format.Range("D6").Select()

# Re-activate Processes so it ends up as the workbook's active sheet/tab.
$processes.Activate()
$processes.Range("D4").Select()

"done"
